$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-30 Tuesday" "2025-10-01 Wednesday"

Replace-Text "474×4=1896" "829×3=2487"
Replace-Text "235×3=705" "270×9=2430"
Replace-Text "589×7=4123" "700×7=4900"
Replace-Text "943×5=4715" "196×4=784"
Replace-Text "135×9=1215" "264×2=528"

Replace-Text "154×3=462" "993×3=2979"
Replace-Text "479×4=1916" "130×6=780"
Replace-Text "910×6=5460" "326×8=2608"
Replace-Text "876×9=7884" "251×7=1757"
Replace-Text "628×5=3140" "273×4=1092"

Replace-Text "485×6=2910" "900×7=6300"
Replace-Text "893×4=3572" "284×7=1988"
Replace-Text "299×6=1794" "677×4=2708"
Replace-Text "536×6=3216" "748×3=2244"
Replace-Text "250×4=1000" "850×7=5950"

Replace-Text "527×3=1581" "424×2=848"
Replace-Text "340×3=1020" "782×5=3910"
Replace-Text "499×5=2495" "155×3=465"
Replace-Text "178×8=1424" "120×7=840"
Replace-Text "300×2=600" "396×3=1188"

Replace-Text "157×2=314" "911×7=6377"
Replace-Text "392×5=1960" "140×2=280"
Replace-Text "971×3=2913" "861×5=4305"
Replace-Text "155×7=1085" "268×8=2144"
Replace-Text "336×8=2688" "325×6=1950"
